# Add a new login-test worksheet ("sauceLoginsTest") between the existing
# "sauceLoginTest" and "InformationPageTest" sheets, seeded with the extra
# SauceDemo test accounts (locked_out_user, problem_user,
# performance_glitch_user) alongside the original standard_user row.

$wb = $excel.ActiveWorkbook

$sauceLoginTest = $wb.Worksheets.Item("sauceLoginTest")

# Inserting with an explicit "After" target places the new sheet right
# after sauceLoginTest (and Excel activates it, which also matches the
# target's activeTab/tabSelected move onto the new sheet).
$newSheet = $wb.Worksheets.Add($null, $sauceLoginTest)
$newSheet.Name = "sauceLoginsTest"

$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "Password"

$newSheet.Range("A2").Value = "standard_user"
$newSheet.Range("B2").Value = "secret_sauce"

$newSheet.Range("A3").Value = "locked_out_user"
$newSheet.Range("B3").Value = "secret_sauce"

$newSheet.Range("A4").Value = "problem_user"
$newSheet.Range("B4").Value = "secret_sauce"

$newSheet.Range("A5").Value = "performance_glitch_user"
$newSheet.Range("B5").Value = "secret_sauce"

# Match the bestFit-ish column widths of the source sheet as closely as
# this host's width quantization allows.
$newSheet.Columns.Item(1).ColumnWidth = 22.6
$newSheet.Columns.Item(2).ColumnWidth = 11.6
